$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 03:08:40"
$wsZhCn.Range("H2").Value = "2016-03-23 03:09:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 03:08:44"
$wsDeDe.Range("H2").Value = "2016-03-23 03:09:16"
